# US-2.3_Accountant_Account payable functionality_TCs.xlsx
# Commit: "Add files via upload" - updates the expected-result text for the
# "vendor dropdown field" validation rows (F19:F21, F46:F48, F68:F70) from
# "It should be displayed..." / the long field-validation sentence to the
# new wording "It should not be displayed and should be as per parameter.",
# and moves the active view/selection down to the newly edited rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "It should not be displayed and should be as per parameter."

# Update the three blocks of "Validate vendor dropdown field" expected results.
$ws.Range("F19").Value = $newText
$ws.Range("F20").Value = $newText
$ws.Range("F21").Value = $newText

$ws.Range("F46").Value = $newText
$ws.Range("F47").Value = $newText
$ws.Range("F48").Value = $newText

$ws.Range("F68").Value = $newText
$ws.Range("F69").Value = $newText
$ws.Range("F70").Value = $newText

# The text in F68:F70 is shorter than the old wrapped text, so the rows
# shrink back down to the default auto-fit height (rows 67 and 71, which
# were not touched, keep their original taller height).
$ws.Rows.Item(68).AutoFit()
$ws.Rows.Item(69).AutoFit()
$ws.Rows.Item(70).AutoFit()

# Scroll the view to the edited area and select the updated cells, matching
# the author's final on-screen position when the file was saved.
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F68:F70").Select()
